# ------------------------------------------------------------------
# Applies the "Built site for gh-pages" regeneration diff:
#   1. Collapse the split title / author / abstract runs into single
#      runs (text content unchanged).
#   2. Bump the lecture date 2024-11-07 -> 2024-11-08.
#   3. Re-point the "Subtitle" style onto "Title" and drop its direct
#      color override (now inherited/automatic).
#   4. Drop the hard-coded blue color on the "AbstractTitle" style.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1a. Title paragraph: merge the word-by-word runs -------------
$d.Paragraphs(1).Range.Find.Execute(
    "Lecture Data Science for Electron Microscopy Winter 2024",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lecture Data Science for Electron Microscopy Winter 2024", 2)

# --- 1b. Author paragraph: merge "Philipp" / " " / "Pelz" ----------
$d.Paragraphs(2).Range.Find.Execute(
    "Philipp Pelz",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Philipp Pelz", 2)

# --- 2. Date paragraph: 2024-11-07 -> 2024-11-08 --------------------
$d.Paragraphs(3).Range.Find.Execute(
    "2024-11-07",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2024-11-08", 2)

# --- 1c. Abstract paragraph: merge the word-by-word runs -----------
$d.Paragraphs(5).Range.Find.Execute(
    "This is the website for the Data Science for Electron Microscopy Lecture",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This is the website for the Data Science for Electron Microscopy Lecture", 2)

# --- 3. Subtitle style: base it on "Title", drop the direct color --
$d.Styles("Subtitle").BaseStyle = "Title"
$d.Styles("Subtitle").Font.ColorIndex = 0

# --- 4. AbstractTitle style: drop the hard-coded 345A8A color -------
$d.Styles("AbstractTitle").Font.ColorIndex = 0
